# edit.ps1 - Apply updated crypto price / 1h-volume data from the
# scheduled GitHub Actions refresh (Fri Feb 16 11:45:25 UTC 2024).
# Includes a handful of ranking swaps (rows 23/24, 29/30/31, 43/44)
# where two/three coins changed relative rank position.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as literal TEXT (preserve the sheet's existing
# inline-string formatting, e.g. "23.30" / "0.999" / "52.349.77") without
# letting Excel auto-coerce number-looking strings into numeric values,
# and without leaving the cell's style/number-format changed afterwards.
function Set-TextValue {
    param($cellRef, $text)
    $r = $ws.Range($cellRef)
    $origStyle = $r.Style
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.Style = $origStyle
}

# --- Price column (D) updates: force-write as text so values like
# "23.30", "0.999", "52.349.77" keep their exact printed form ---
$priceUpdates = [ordered]@{
    "D2" = "52.349.77"
    "D3" = "2.831.99"
    "D5" = "356.36"
    "D6" = "112.54"
    "D7" = "0.573"
    "D8" = "0.999"
    "D9" = "0.601"
    "D10" = "41.15"
    "D13" = "19.96"
    "D14" = "7.79"
    "D15" = "3.274.04"
    "D16" = "2.836.45"
    "D18" = "52.154.82"
    "D19" = "7.52"
    "D21" = "13.53"
    "D22" = "0.0₃0998"
    "D23" = "70.66"
    "D24" = "271.82"
    "D25" = "2.81"
    "D26" = "27.01"
    "D28" = "10.36"
    "D29" = "2.26"
    "D30" = "0.0491"
    "D31" = "0.144"
    "D32" = "52.61"
    "D33" = "35.12"
    "D34" = "5.95"
    "D36" = "0.0858"
    "D38" = "3.27"
    "D40" = "18.46"
    "D41" = "0.118"
    "D42" = "127.42"
    "D43" = "2.54"
    "D44" = "23.30"
    "D45" = "2.28"
    "D46" = "3.37"
    "D47" = "2.090.62"
    "D48" = "2.28"
    "D49" = "5.95"
    "D50" = "0.974"
}
foreach ($cell in $priceUpdates.Keys) {
    Set-TextValue $cell $priceUpdates[$cell]
}

# --- Coin name (B), link (C) and 1h-volume % (E) updates ---
$otherUpdates = [ordered]@{
    "E2" = "  +0.18%  "
    "E3" = "  +1.51%  "
    "E4" = "  -0.05%  "
    "E5" = "  +2.96%  "
    "E6" = "  -3.05%  "
    "E7" = "  +4.19%  "
    "E8" = "  -0.01%  "
    "E9" = "  +1.86%  "
    "E10" = "  -3.93%  "
    "E11" = "  +1.00%  "
    "E12" = "  +1.15%  "
    "E13" = "  -0.41%  "
    "E14" = "  -0.63%  "
    "E15" = "  +1.42%  "
    "E16" = "  +1.19%  "
    "E17" = "  +4.34%  "
    "E18" = "  +0.06%  "
    "E19" = "  +4.54%  "
    "E20" = "  -0.45%  "
    "E21" = "  +0.95%  "
    "E22" = "  +1.88%  "
    "B23" = "Litecoin"
    "C23" = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
    "E23" = "  +0.96%  "
    "B24" = "BitcoinCash"
    "C24" = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
    "E24" = "  +0.93%  "
    "E25" = "  +1.82%  "
    "E26" = "  +1.45%  "
    "E27" = "  -0.05%  "
    "E28" = "  +1.68%  "
    "B29" = "Toncoin"
    "C29" = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
    "E29" = "  +0.55%  "
    "B30" = "VeChain"
    "C30" = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
    "E30" = "  +9.93%  "
    "B31" = "Kaspa"
    "C31" = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
    "E31" = "  +3.02%  "
    "E32" = "  +4.77%  "
    "E33" = "  +0.63%  "
    "E34" = "  +4.31%  "
    "E35" = "  +12.60%  "
    "E36" = "  +4.03%  "
    "E37" = "  -0.21%  "
    "E38" = "  +1.84%  "
    "E39" = "  -3.05%  "
    "E40" = "  -1.36%  "
    "E41" = "  +2.22%  "
    "E42" = "  -0.52%  "
    "B43" = "Stacks"
    "C43" = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
    "E43" = "  -4.98%  "
    "B44" = "EnergySwap"
    "C44" = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
    "E44" = "  -1.82%  "
    "E45" = "  -1.31%  "
    "E46" = "  +1.20%  "
    "E47" = "  +1.23%  "
    "E48" = "  -2.62%  "
    "E49" = "  +8.04%  "
    "E50" = "  +0.51%  "
    "E51" = "  +3.03%  "
}
foreach ($cell in $otherUpdates.Keys) {
    $ws.Range($cell).Value = $otherUpdates[$cell]
}
